$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the used range and update any cell in column B whose text is
# "Land & Water Conservation Fund" to the new label "Land & Water Conservation".
$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "Land & Water Conservation Fund") {
        $cell.Value = "Land & Water Conservation"
    }
}

# Widen column B to fit the (differently sized) label text, matching the
# authored width change (23.109375 -> 33.6640625 raw OOXML width units).
# ColumnWidth is specified in characters and gets quantized internally to
# whole pixels, so 32.83 is the closest achievable value that reproduces
# the target raw width.
$ws.Columns.Item(2).ColumnWidth = 32.83
